$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds "line1".."line6" (rows 2-7) followed directly by
# "extr1".."extr8" (rows 8-15). Two new rows for "line7" and "line8" need to
# be inserted between them, pushing extr1..extr8 down by two rows (to rows
# 10-17) and updating their data.

# Shift rows 8-15 (extr1..extr8) down to rows 10-17 by copying whole rows,
# working from the bottom up so data is not overwritten before it is copied.
for ($r = 15; $r -ge 8; $r--) {
    $dest = $r + 2
    $srcRange = $ws.Range("A$r" + ":E$r")
    $destRange = $ws.Range("A$dest" + ":E$dest")
    $srcRange.Copy($destRange)
}

# Renumber the sequential index column (A) for the shifted extr1..extr8 rows
# (now at rows 10-17) so it keeps counting up from the new line7/line8 rows.
$newIndex = @(8, 9, 10, 11, 12, 13, 14, 15)
for ($i = 0; $i -lt 8; $i++) {
    $row = 10 + $i
    $ws.Cells.Item($row, 1).Value = $newIndex[$i]
}

# Fill in the two new rows: line7 (row 8) and line8 (row 9).
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# Update the in_service flag (column E) for the shifted extr1..extr8 rows.
$ws.Cells.Item(10, 5).Value = $true    # extr1
$ws.Cells.Item(11, 5).Value = $true    # extr2
$ws.Cells.Item(12, 5).Value = $false   # extr3
$ws.Cells.Item(13, 5).Value = $false   # extr4
$ws.Cells.Item(14, 5).Value = $false   # extr5
$ws.Cells.Item(15, 5).Value = $true    # extr6
$ws.Cells.Item(16, 5).Value = $false   # extr7
$ws.Cells.Item(17, 5).Value = $false   # extr8
